# Update the Volume(1h) percentage values in column E of the cryptos sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "  -2.27%  "
$ws.Range("E3").Value = "  -5.01%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("E5").Value = "  -1.45%  "
$ws.Range("E6").Value = "  -7.55%  "
$ws.Range("E8").Value = "  -10.83%  "
$ws.Range("E9").Value = "  -4.98%  "
$ws.Range("E10").Value = "  -2.25%  "
$ws.Range("E11").Value = "  -0.10%  "
$ws.Range("E12").Value = "  -2.47%  "
$ws.Range("E13").Value = "  -3.63%  "
$ws.Range("E14").Value = "  -6.73%  "
$ws.Range("E15").Value = "  -5.02%  "
$ws.Range("E16").Value = "  -1.69%  "
$ws.Range("E17").Value = "  -4.56%  "
$ws.Range("E18").Value = "  -4.91%  "
$ws.Range("E19").Value = "  -4.11%  "
$ws.Range("E20").Value = "  -1.23%  "
$ws.Range("E21").Value = "  -2.79%  "
$ws.Range("E22").Value = "  -6.63%  "
$ws.Range("E23").Value = "  -0.11%  "
$ws.Range("E24").Value = "  -2.85%  "
$ws.Range("E25").Value = "  -1.59%  "
$ws.Range("E26").Value = "  +10.75%  "
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("E28").Value = "  -5.07%  "
$ws.Range("E29").Value = "  -9.46%  "
$ws.Range("E30").Value = "  -4.02%  "
$ws.Range("E31").Value = "  -6.38%  "
$ws.Range("E32").Value = "  -6.53%  "
$ws.Range("E34").Value = "  -5.13%  "
$ws.Range("E37").Value = "  -3.91%  "
$ws.Range("E38").Value = "  -2.20%  "
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("E40").Value = "  -9.60%  "
$ws.Range("E41").Value = "  +1.89%  "
$ws.Range("E42").Value = "  -0.08%  "
$ws.Range("E43").Value = "  -2.11%  "
$ws.Range("E44").Value = "  -0.82%  "
$ws.Range("E45").Value = "  -2.04%  "
$ws.Range("E46").Value = "  -7.40%  "
$ws.Range("E47").Value = "  -2.58%  "
$ws.Range("E48").Value = "  -4.48%  "
$ws.Range("E49").Value = "  -10.98%  "
$ws.Range("E50").Value = "  -4.12%  "
$ws.Range("E51").Value = "  -3.72%  "
